$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("J3").Value = 2.25
$ws.Range("K3").Value = 2.5
$ws.Range("L3").Value = 4.33
$ws.Range("S3").Value = 1.25
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 9
$ws.Range("AP3").Value = 15
$ws.Range("AQ3").Value = 26
$ws.Range("AR3").Value = 41
$ws.Range("AS3").Value = 81
$ws.Range("AT3").Value = 3.75
$ws.Range("AU3").Value = 7
$ws.Range("AV3").Value = 41
$ws.Range("AW3").Value = 6.5
$ws.Range("AX3").Value = 21
$ws.Range("AY3").Value = 23
$ws.Range("AZ3").Value = 51
$ws.Range("BA3").Value = 67
$ws.Range("BB3").Value = 126
$ws.Range("BC3").Value = 301

# Row 6
$ws.Range("N6").Value = 8
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 1.62

# Row 37
$ws.Range("N37").Value = 12
$ws.Range("U37").Value = 1.75
$ws.Range("V37").Value = 2
$ws.Range("AC37").Value = 12
$ws.Range("AO37").Value = 23
$ws.Range("AX37").Value = 9

# Row 38
$ws.Range("G38").Value = 1.73
$ws.Range("J38").Value = 2.3
$ws.Range("L38").Value = 4.75
$ws.Range("M38").Value = 1.05
$ws.Range("N38").Value = 11
$ws.Range("Q38").Value = 1.85
$ws.Range("R38").Value = 2
$ws.Range("AB38").Value = 23
$ws.Range("AC38").Value = 11
$ws.Range("AD38").Value = 7
$ws.Range("AG38").Value = 201
$ws.Range("AI38").Value = 23

# Row 39
$ws.Range("M39").Value = 8.3

# Row 49
$ws.Range("G49").Value = 1.67
$ws.Range("I49").Value = 5.25
$ws.Range("J49").Value = 2.25
$ws.Range("K49").Value = 2.25
$ws.Range("L49").Value = 5
$ws.Range("Q49").Value = 1.85
$ws.Range("R49").Value = 2
$ws.Range("Z49").Value = 13
$ws.Range("AG49").Value = 251
$ws.Range("AR49").Value = 51
$ws.Range("AW49").Value = 6.5
$ws.Range("AZ49").Value = 81
$ws.Range("BC49").Value = 126
